# The "Mapp_Differences" sheet lists records that could not be matched
# between mapp_api.csv and mapp_prts.csv. Rows 136-142 are the block of
# "ID not found in mapp_api.csv" records; this edit reorders that block
# (columns A, C, D and E - the ID, the raw record, the created_At epoch,
# and the created_At date/time). Columns B ("N/A") and F
# ("ID not found in mapp_api.csv") are identical for every row in the
# block, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $id, $record, $epoch, $dateTime) {
    $ws.Range("A$row").Value = $id
    $ws.Range("C$row").Value = $record

    # $epoch is a digit-only string ("1730798554", ...) - a bare .Value
    # assignment would make Excel re-interpret it as a number, losing the
    # shared-string text type. Force text via NumberFormat, then clear the
    # formatting again so the cell's style stays the default (matches the
    # rest of the sheet, which carries no explicit style).
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $epoch
    $ws.Range("D$row").ClearFormats()

    $ws.Range("E$row").Value = $dateTime
}

Set-Row 136 "PHZYMuYZy81kiB" "[PHZYMuYZy81kiB PHTl2JAhUSu8s6 merchant PHZYM0x70Hkl2q 1730798554 1730798554 ]" "1730798554" "2024-11-05 14:52:34"
Set-Row 137 "P4YgZMYhoO9C6a" "[P4YgZMYhoO9C6a C65yyS9pcolRUr referred P4YgYRvTyeMzWs 1727957076 1727957076 ]" "1727957076" "2024-10-03 17:34:36"
Set-Row 138 "P62kxToRP7l6xm" "[P62kxToRP7l6xm C65yyS9pcolRUr referred P62kwFoO4o2awa 1728281314 1728281314 ]" "1728281314" "2024-10-07 11:38:34"
Set-Row 139 "P4YjmZK5ukGRTR" "[P4YjmZK5ukGRTR C65yyS9pcolRUr referred P4YjleeSMz6H65 1727957259 1727957259 ]" "1727957259" "2024-10-03 17:37:39"
Set-Row 140 "PHZXGEvIsF1zmx" "[PHZXGEvIsF1zmx PHWlZROPA0es90 merchant PHZXFL1AQ7fp1q 1730798491 1730798491 ]" "1730798491" "2024-11-05 14:51:31"
Set-Row 141 "PHZSWu73Wn9FkZ" "[PHZSWu73Wn9FkZ PHZRkFuWQ2kKhs referred PHZSUvYf6ngG3k 1730798222 1730798222 ]" "1730798222" "2024-11-05 14:47:02"
Set-Row 142 "P4rb4ZA28eywM3" "[P4rb4ZA28eywM3 C65yyS9pcolRUr referred P4rb3LbPo3hd86 1728023675 1728023675 ]" "1728023675" "2024-10-04 12:04:35"
